$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "conditioned", 434, 0, 1, "KM"),
    @(1, "unconditioned", 434, 0, 0, "EF"),
    @(2, "conditioned", 434, 0, 1, "KM"),
    @(2, "unconditioned", 434, 0, 0, "EF"),
    @(3, "conditioned", 434, 0, 0, "KM"),
    @(3, "unconditioned", 434, 0, 0, "EF"),
    @(4, "conditioned", 434, 0, 1, "KM"),
    @(4, "unconditioned", 434, 0, 0, "EF"),
    @(5, "conditioned", 434, 0, 0, "KM"),
    @(5, "unconditioned", 434, 0, 0, "EF"),
    @(6, "conditioned", 434, 0, 0, "EF"),
    @(6, "unconditioned", 434, 0, 0, "KM"),
    @(7, "conditioned", 434, 0, 0, "EF"),
    @(7, "unconditioned", 434, 0, 0, "KM"),
    @(8, "conditioned", 434, 0, 0, "EF"),
    @(8, "unconditioned", 434, 0, 0, "KM"),
    @(9, "conditioned", 434, 0, 2, "EF"),
    @(9, "unconditioned", 434, 0, 0, "KM"),
    @(10, "conditioned", 434, 0, 0, "KM"),
    @(10, "unconditioned", 434, 0, 0, "KM"),
    @(11, "conditioned", 434, 0, 0, "KM"),
    @(11, "unconditioned", 434, 0, 1, "EF"),
    @(12, "conditioned", 434, 0, 0, "EF"),
    @(12, "unconditioned", 434, 0, 0, "EF"),
    @(13, "conditioned", 434, 0, 0, "EF"),
    @(13, "unconditioned", 434, 0, 0, "KM"),
    @(14, "conditioned", 434, 0, 0, "EF"),
    @(14, "unconditioned", 434, 0, 0, "KM"),
    @(15, "unconditioned", 434, 0, 0, "KM")
)

$startRow = 379
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
}

$excel.ActiveWindow.Zoom = 194
$excel.ActiveWindow.ScrollRow = 393
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F403").Select()
